# Update header labels on the existing sheets
$wb = $excel.ActiveWorkbook

$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# Add the new "PO Forecast" sheet after "Monthly Trend" (last sheet)
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count())
$wsForecast = $wb.Worksheets.Add($null, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Header row
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Match the bold/centered/bordered header style already used on the other sheets
$wsWeekly.Range("A1:B1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)

# Forecast data rows (ds, PO_Forecast, yhat_lower, yhat_upper)
$data = @(
    @(44934.99999999999, 121, -82.25860060225422, 309.6237101124069),
    @(44948.99999999999, 121, -78.63619191129482, 337.9991507008518),
    @(44969.99999999999, 122, -70.55914297052072, 340.9030775801112),
    @(44976.99999999999, 122, -79.65618469095149, 338.7536816161659),
    @(44983.99999999999, 122, -90.7504312270907, 328.7256828625875),
    @(44990.99999999999, 122, -76.6377571858342, 308.898420708767),
    @(45011.99999999999, 122, -85.09233733881503, 330.2614796803527),
    @(45018.99999999999, 122, -84.00809551558042, 317.7692243151095),
    @(45025.99999999999, 122, -84.48269973072262, 331.9520558468052),
    @(45032.99999999999, 123, -84.23076761233808, 333.6341963773268),
    @(45053.99999999999, 123, -80.60773779590696, 330.9653484035785),
    @(45060.99999999999, 123, -77.86822596581617, 329.1454590421457),
    @(45067.99999999999, 123, -86.93127399171507, 334.2928859868883),
    @(45074.99999999999, 123, -80.42068219566839, 333.5770490627845),
    @(45081.99999999999, 123, -91.38663106198798, 320.9082214922648),
    @(45088.99999999999, 123, -71.32210739565114, 342.3568192507882),
    @(45095.99999999999, 124, -85.100955239779, 329.27206413739),
    @(45102.99999999999, 124, -77.25671327635919, 316.5255343185956),
    @(45109.99999999999, 124, -66.37546689425754, 331.6494675746382),
    @(45116.99999999999, 124, -100.168398917165, 327.4111884375754),
    @(45123.99999999999, 124, -73.781689187683, 331.1629970382531),
    @(45130.99999999999, 124, -99.10098615883284, 334.8243841673545),
    @(45137.99999999999, 124, -81.99910008753322, 311.6499313018745),
    @(45144.99999999999, 124, -84.03627612289866, 325.0505360931293),
    @(45151.99999999999, 124, -82.83470334591347, 324.3142765272408),
    @(45158.99999999999, 125, -70.75962890981174, 322.5945437153524),
    @(45165.99999999999, 125, -70.22473049047416, 334.8373826599922),
    @(45172.99999999999, 125, -70.67505929929037, 331.0029132701269),
    @(45179.99999999999, 125, -71.67550149772993, 323.5878417157383),
    @(45186.99999999999, 125, -62.06500373085046, 330.394384075586),
    @(45193.99999999999, 125, -89.64038701354073, 326.576604211434),
    @(45207.99999999999, 125, -70.57088901558346, 335.6086281264288),
    @(45214.99999999999, 125, -93.73672975140985, 323.588360181679),
    @(45221.99999999999, 126, -82.88957602923647, 337.1958900331504),
    @(45228.99999999999, 126, -77.83510227385135, 318.8768758237147),
    @(45235.99999999999, 126, -86.51203049669762, 336.4253339470947),
    @(45242.99999999999, 126, -80.51564576941709, 332.7457766900451),
    @(45249.99999999999, 126, -79.88942685875753, 328.6290165099157),
    @(45256.99999999999, 126, -89.18831885585861, 329.1217619426961),
    @(45263.99999999999, 126, -66.71037493485537, 325.9134617245405),
    @(45270.99999999999, 126, -78.65722130363392, 326.8584711804953),
    @(45277.99999999999, 126, -86.61528672178963, 319.6498885222162),
    @(45298.99999999999, 127, -78.53015311648589, 352.0809927331902),
    @(45305.99999999999, 127, -80.4786831843314, 341.9493747589128),
    @(45312.99999999999, 127, -73.8652337443404, 333.3700495481942),
    @(45319.99999999999, 127, -69.2203168237464, 317.5452429782047),
    @(45326.99999999999, 127, -73.52139952982034, 329.1049139198244),
    @(45333.99999999999, 127, -77.98054104487638, 332.3754729016367),
    @(45340.99999999999, 127, -80.46745600779863, 312.0409321309349),
    @(45347.99999999999, 128, -80.20095983520116, 331.7348304003558),
    @(45354.99999999999, 128, -55.44546389048069, 332.0535247798717),
    @(45361.99999999999, 128, -72.5098931458089, 344.697971190452),
    @(45368.99999999999, 128, -71.8599565675343, 332.5470227384062),
    @(45375.99999999999, 128, -72.45135181790474, 344.900050928595),
    @(45382.99999999999, 128, -67.93034494564056, 333.8094440135945),
    @(45389.99999999999, 128, -80.46053016574059, 335.5362891454716),
    @(45417.99999999999, 129, -80.97584424994636, 347.9046540305912),
    @(45431.99999999999, 129, -68.80481909482721, 334.8298250771123),
    @(45445.99999999999, 129, -79.60158821050186, 326.2721743730879),
    @(45452.99999999999, 129, -82.08773054314992, 315.1368441743171),
    @(45487.99999999999, 130, -73.74056746172026, 330.1922175411883),
    @(45494.99999999999, 130, -71.7416044942026, 331.1482543207096),
    @(45536.99999999999, 131, -75.5914794368291, 347.7031781378753),
    @(45543.99999999999, 131, -85.50911493183298, 316.4291587953825),
    @(45578.99999999999, 131, -79.45133116000653, 333.198989940358),
    @(45592.99999999999, 131, -80.22848586383303, 325.7517996782803),
    @(45613.99999999999, 132, -78.31814646953062, 333.7220785076565),
    @(45620.99999999999, 132, -66.97966407509634, 331.431636322415),
    @(45641.99999999999, 132, -74.92580342075803, 343.4569693402465),
    @(45648.99999999999, 132, -71.86066277546591, 328.2583151217169),
    @(45655.99999999999, 132, -80.22214052001618, 331.6967473341119),
    @(45662.99999999999, 133, -91.74538417235561, 333.8624743414553),
    @(45669.99999999999, 133, -61.95201817607147, 350.2895913552278),
    @(45676.99999999999, 133, -75.93239766120357, 345.030053600871),
    @(45683.99999999999, 133, -62.11958287274546, 344.7423141382511),
    @(45690.99999999999, 133, -75.48338245968047, 337.5256286870272),
    @(45697.99999999999, 133, -92.75224825829805, 342.4323282863716)
)

$row = 2
foreach ($r in $data) {
    $wsForecast.Cells.Item($row, 1).Value = $r[0]
    $wsForecast.Cells.Item($row, 2).Value = $r[1]
    $wsForecast.Cells.Item($row, 3).Value = $r[2]
    $wsForecast.Cells.Item($row, 4).Value = $r[3]
    $row = $row + 1
}

$dateRange = $wsForecast.Range("A2:A78")
$dateRange.NumberFormat = "YYYY-MM-DD HH:MM:SS"
